# Updating cart lang and fb pixel test data and updating buyflow output
# sheet including Status column.
#
# Sheet1: the second (duplicate) results table that used to live at rows
# 24-39 is removed; four fresh rows describing the new
# "deluxe25offp-redes"/"cpcb2017"/"cpwbunusedbdbj" kit test cases are
# inserted right after the existing QA/Sub-D rows (12-14), just before the
# trailing "End" marker.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Drop the blank spacer row (16) plus the old duplicate table (rows 17-39),
# leaving just the "End" marker row (still row 15) behind.
$ws1.Range("A16:A39").EntireRow.Delete()

# Make room for four new data rows ahead of the "End" marker, shifting it
# from row 15 down to row 19.
$ws1.Range("A15:E15").EntireRow.Insert()
$ws1.Range("A15:E15").EntireRow.Insert()
$ws1.Range("A15:E15").EntireRow.Insert()
$ws1.Range("A15:E15").EntireRow.Insert()

$ws1.Cells.Item(15,1).Value = "QA"
$ws1.Cells.Item(15,2).Value = "Sub-D"
$ws1.Cells.Item(15,3).Value = "cpcb2017"
$ws1.Cells.Item(15,4).Value = "Kit"
$ws1.Cells.Item(15,5).Value = "Chrome"

$ws1.Cells.Item(16,1).Value = "QA"
$ws1.Cells.Item(16,2).Value = "Sub-D"
$ws1.Cells.Item(16,3).Value = "deluxe25offp-redes"
$ws1.Cells.Item(16,4).Value = "Kit"
$ws1.Cells.Item(16,5).Value = "Chrome"

$ws1.Cells.Item(17,1).Value = "QA"
$ws1.Cells.Item(17,2).Value = "Sub-D"
$ws1.Cells.Item(17,3).Value = "deluxe25offp"
$ws1.Cells.Item(17,4).Value = "Kit"
$ws1.Cells.Item(17,5).Value = "Chrome"

$ws1.Cells.Item(18,1).Value = "QA"
$ws1.Cells.Item(18,2).Value = "Sub-D"
$ws1.Cells.Item(18,3).Value = "cpwbunusedbdbj"
$ws1.Cells.Item(18,4).Value = "Kit"
$ws1.Cells.Item(18,5).Value = "Chrome"

# Move the selection to reflect where editing left off.
$ws1.Range("B17").Select() | Out-Null

# R2 sheet: selection now highlights the entire 6th row instead of the old
# A22:E22 block.
$ws2 = $wb.Worksheets.Item(3)
$ws2.Rows.Item(6).Select() | Out-Null

# Leave the workbook focused back on Sheet1 (the tab that was active
# before/after the edit).
$ws1.Activate() | Out-Null
